$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 842
$ws.Cells.Item(3, 6).Value = 1737
$ws.Cells.Item(4, 6).Value = 35
$ws.Cells.Item(5, 6).Value = 529
$ws.Cells.Item(6, 6).Value = 2144
$ws.Cells.Item(7, 6).Value = 1358
$ws.Cells.Item(8, 6).Value = 2058
$ws.Cells.Item(9, 6).Value = 962
$ws.Cells.Item(11, 6).Value = 2394
$ws.Cells.Item(12, 6).Value = 653
$ws.Cells.Item(14, 6).Value = 3867
$ws.Cells.Item(16, 6).Value = 362
$ws.Cells.Item(17, 6).Value = 2934
$ws.Cells.Item(18, 6).Value = 773
$ws.Cells.Item(19, 6).Value = 139
$ws.Cells.Item(21, 6).Value = 96
$ws.Cells.Item(22, 6).Value = 2015
$ws.Cells.Item(23, 6).Value = 1166
$ws.Cells.Item(24, 6).Value = 1815
$ws.Cells.Item(25, 6).Value = 376
$ws.Cells.Item(27, 6).Value = 5
$ws.Cells.Item(28, 6).Value = 8191
$ws.Cells.Item(29, 6).Value = 5571
$ws.Cells.Item(30, 6).Value = 348
$ws.Cells.Item(31, 6).Value = 164
$ws.Cells.Item(32, 6).Value = 744
$ws.Cells.Item(33, 6).Value = 756
$ws.Cells.Item(36, 6).Value = 942
$ws.Cells.Item(37, 6).Value = 377
$ws.Cells.Item(38, 6).Value = 27
$ws.Cells.Item(39, 6).Value = 184
$ws.Cells.Item(40, 6).Value = 147
$ws.Cells.Item(41, 6).Value = 4596
$ws.Cells.Item(42, 6).Value = 823
$ws.Cells.Item(43, 6).Value = 61

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(3, 6).Value = 86
$ws.Cells.Item(18, 6).Value = 161
$ws.Cells.Item(27, 6).Value = 5

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 8213
$ws.Cells.Item(4, 6).Value = 1244

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 8213
$ws.Cells.Item(3, 6).Value = 842
$ws.Cells.Item(5, 6).Value = 1244
$ws.Cells.Item(6, 6).Value = 86
$ws.Cells.Item(7, 6).Value = 1737
$ws.Cells.Item(8, 6).Value = 35
$ws.Cells.Item(9, 6).Value = 529
$ws.Cells.Item(10, 6).Value = 1358
$ws.Cells.Item(11, 6).Value = 2058
$ws.Cells.Item(12, 6).Value = 962
$ws.Cells.Item(15, 6).Value = 3867
$ws.Cells.Item(16, 6).Value = 362
$ws.Cells.Item(17, 6).Value = 2934
$ws.Cells.Item(18, 6).Value = 773
$ws.Cells.Item(20, 6).Value = 2015
$ws.Cells.Item(26, 6).Value = 1166
$ws.Cells.Item(28, 6).Value = 1815
$ws.Cells.Item(30, 6).Value = 376
$ws.Cells.Item(31, 6).Value = 8191
$ws.Cells.Item(32, 6).Value = 5571
$ws.Cells.Item(34, 6).Value = 348
$ws.Cells.Item(35, 6).Value = 164
$ws.Cells.Item(36, 6).Value = 744
$ws.Cells.Item(37, 6).Value = 756
$ws.Cells.Item(39, 6).Value = 942
$ws.Cells.Item(40, 6).Value = 377
$ws.Cells.Item(41, 6).Value = 184
$ws.Cells.Item(43, 6).Value = 4596
$ws.Cells.Item(44, 6).Value = 823
$ws.Cells.Item(45, 6).Value = 61
$ws.Cells.Item(49, 6).Value = 5
